# edit.ps1 - applies the "Added new project and updated resume" changes
# to brian-h-resume.docx via the Word COM-interop object model.

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, [string]$needle) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) Consultant Web Developer job end date: "-Present" -> "-March 2017"
#    (scoped to the specific paragraph, since "-Present" also occurs in
#     the EDUCATION "2016-Present" line)
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "March 2015"
$p = $d.Paragraphs.Item($idx)
$null = $p.Range.Find.Execute("-Present", $true, $false, $false, $false, $false, $true, 1, $false, "-March 2017", 2)

# ---------------------------------------------------------------------
# 2) "...jQuery, HTML, and SQL." -> "...jQuery, HTML5, and SQL."
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(", jQuery, HTML, and SQL.", $true, $false, $false, $false, $false, $true, 1, $false, ", jQuery, HTML5, and SQL.", 2)

# ---------------------------------------------------------------------
# 3) "...using Node.js and React.js." -> "...using Node.js, React.js, Socket.IO, and RabbitMQ."
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute("real-time data display, using Node.js and React.js.", $true, $false, $false, $false, $false, $true, 1, $false, "real-time data display, using Node.js, React.js, Socket.IO, and RabbitMQ.", 2)

# ---------------------------------------------------------------------
# 4) "... library), JSON, and SQL." -> "... library), JSON, SQL, and Bootstrap."
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(" library), JSON, and SQL.", $true, $false, $false, $false, $false, $true, 1, $false, " library), JSON, SQL, and Bootstrap.", 2)

# ---------------------------------------------------------------------
# 5) Delete the whole bullet paragraph:
#    "Supported users in a fast-paced trading environment."
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "Supported users in a fast-paced trading environment."
$p = $d.Paragraphs.Item($idx)
$p.Range.Delete()

# ---------------------------------------------------------------------
# 6) "...jQuery, and JSON." -> "...jQuery, JSON, and XML."
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(", JavaScript, jQuery, and JSON.", $true, $false, $false, $false, $false, $true, 1, $false, ", JavaScript, jQuery, JSON, and XML.", 2)

# ---------------------------------------------------------------------
# 7) "full-stack applications with Node.js, React.js, and PostgreSQL." ->
#    "full-stack applications with Node.js, Express.js, React.js, PostgreSQL, and EJS."
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute("full-stack applications with Node.js, React.js, and PostgreSQL.", $true, $false, $false, $false, $false, $true, 1, $false, "full-stack applications with Node.js, Express.js, React.js, PostgreSQL, and EJS.", 2)

# ---------------------------------------------------------------------
# 8) Insert a new PROJECTS bullet before the "small classroom library" one:
#    "Node.js, EJS, and PostgreSQL application for encoding poetry in TEI. (CUNY - Available on GitHub)"
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "small classroom library"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphBefore()
$newP = $d.Paragraphs.Item($idx)
$newP.Range.Text = "Node.js, EJS, and PostgreSQL application for encoding poetry in TEI. (CUNY " + [char]0x2013 + " Available on GitHub)"

Write-Host "Edits applied."
